# This script reproduces the target diff:
#   1. Removes the stray "_GoBack" bookmark from its original position
#      (right after the "...שמנוסחים בצורה מסודרת..." run in the 2nd
#      paragraph) since the diff deletes it from there.
#   2. Turns the final (empty) paragraph of the document into a
#      separator line of dashes, and appends a brand-new paragraph after
#      it containing the new Hebrew announcement text, ending with the
#      "_GoBack" bookmark re-inserted at its new location.

$d = $word.ActiveDocument

# --- Step 1: drop the "_GoBack" bookmark from its old spot -----------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# --- Step 2: build the two new paragraphs as a WordOpenXML fragment and ----
#             insert them right before the trailing (empty) paragraph.
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($n)
$insertionPoint = $lastPara.Range

$newContentXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi" w:hint="cs"/><w:rtl/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi" w:hint="cs"/><w:rtl/></w:rPr><w:lastRenderedPageBreak/><w:t>------------------------------------------------------------------------------------------------------------</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi" w:hint="cs"/><w:rtl/></w:rPr><w:t xml:space="preserve">שלום לכולם. </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi"/><w:rtl/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi" w:hint="cs"/><w:rtl/></w:rPr><w:t xml:space="preserve">בשנים האחרונות שמעתי המון תלונות וטענות על מצב התעסוקה באזורנו. כמו רבים מאיתנו, יצא לי ולרוב האנשים שאני מכיר לחפש </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi" w:hint="cs"/><w:rtl/></w:rPr><w:t>עבודה ולא מעט. אף פעם לא קל לחפש עבודה.</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi" w:hint="cs"/><w:rtl/></w:rPr><w:t xml:space="preserve"> כדי לעזור לאנשים עם תעסוקה, פתחתי קבוצת תעסוקה. כיום, בזכות השקעה עצומה אנחנו מגיעים להיקף פעילות יחסית גבוה. בקבוצה שלנו</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$insertionPoint.InsertXML($newContentXml)

# --- Step 3: the original trailing empty paragraph is still hanging ----
#             around after the two freshly-inserted paragraphs; merge it
#             away by deleting its paragraph mark so the new Hebrew
#             paragraph becomes the document's final paragraph (taking
#             over that empty paragraph's own paragraph mark/formatting,
#             exactly like the diff shows).
$newCount = $d.Paragraphs.Count
$hebrewPara = $d.Paragraphs($newCount - 1)
$mergeRange = $d.Range($hebrewPara.Range.End - 1, $hebrewPara.Range.End)
$mergeRange.Delete()
